$wb = $excel.ActiveWorkbook

# Sheet "Metadata": update the ValueSet URL and Date values
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://example.org/ig/example/ValueSet/presence-valueset"
$wsMeta.Range("B8").Value = "2023-04-26T11:15:05-05:00"

# Sheet "Include ValueSets": update the Snomed CT ValueSet URL
$wsSct = $wb.Worksheets.Item("Include ValueSets")
$wsSct.Range("A2").Value = "http://example.org/ig/example/ValueSet/presence-sct-valueset"

# Sheet "Include ValueSets 2": update the LOINC ValueSet URL
$wsLnc = $wb.Worksheets.Item("Include ValueSets 2")
$wsLnc.Range("A2").Value = "http://example.org/ig/example/ValueSet/presence-lnc-valueset"
